$wb = $excel.ActiveWorkbook

# --- Add the new "CreateAccount" worksheet after "SignIn" ---
$wsSignIn = $wb.Worksheets.Item("SignIn")
$wsNew = $wb.Worksheets.Add()
$wsNew.Name = "CreateAccount"
$wsNew.Move($null, $wsSignIn)
$wb.Worksheets.Item("CreateAccount").Move($wsSignIn.Index + 1)

$ws2 = $wb.Worksheets.Item("CreateAccount")

# Header row
$ws2.Range("A1").Value = "RunMode"
$ws2.Range("B1").Value = "TestCase"
$ws2.Range("C1").Value = "FirstName"
$ws2.Range("D1").Value = "LastName"
$ws2.Range("E1").Value = "DateofBirth"
$ws2.Range("F1").Value = "Email"
$ws2.Range("G1").Value = "Company"
$ws2.Range("H1").Value = "Address"
$ws2.Range("I1").Value = "Zipcode"
$ws2.Range("J1").Value = "City"
$ws2.Range("K1").Value = "State"
$ws2.Range("L1").Value = "Country"
$ws2.Range("M1").Value = "TelPhoneNo"
$ws2.Range("N1").Value = "Password"
$ws2.Range("O1").Value = "ConfirmPassword"

# Data row (order matters for shared-string table insertion order)
$ws2.Range("A2").Value = "Y"
$ws2.Range("B2").Value = "CreateAccount"
$ws2.Range("C2").Value = "John"
$ws2.Range("D2").Value = "Smith"
$ws2.Range("G2").Value = "Software Company "
$ws2.Range("H2").Value = "Sholinganallur"
$ws2.Range("J2").Value = "Chennai"
$ws2.Range("K2").Value = "Tamil Nadu"
$ws2.Range("L2").Value = "India"
$ws2.Range("N2").Value = "jaga@12345"
$ws2.Range("E2").NumberFormat = "@"
$ws2.Range("E2").Value = "08/29/1993"
$ws2.Range("F2").Value = "testjaga006@gmail.com"
$ws2.Range("I2").Value = 600119
$ws2.Range("M2").Value = 1234567890
$ws2.Range("O2").Value = "jaga@12345"

# Hyperlinks
$ws2.Hyperlinks.Add($ws2.Range("F2"), "mailto:testjaga006@gmail.com")
$ws2.Hyperlinks.Add($ws2.Range("N2"), "mailto:jaga@12345")
$ws2.Hyperlinks.Add($ws2.Range("O2"), "mailto:jaga@12345")
$ws2.Range("F2").Style = "Hyperlink"
$ws2.Range("N2").Style = "Hyperlink"
$ws2.Range("O2").Style = "Hyperlink"

# Column widths (approximate auto-fit behavior from diff)
$ws2.Range("B1").EntireColumn.ColumnWidth = 14.140625
$ws2.Range("C1").EntireColumn.ColumnWidth = 15.28515625
$ws2.Range("D1").EntireColumn.ColumnWidth = 10.140625
$ws2.Range("E1").EntireColumn.ColumnWidth = 11.28515625
$ws2.Range("F1").EntireColumn.ColumnWidth = 22.42578125
$ws2.Range("G1").EntireColumn.ColumnWidth = 18.42578125
$ws2.Range("H1").EntireColumn.ColumnWidth = 13.85546875
$ws2.Range("K1").EntireColumn.ColumnWidth = 11
$ws2.Range("M1").EntireColumn.ColumnWidth = 12
$ws2.Range("N1").EntireColumn.ColumnWidth = 12
$ws2.Range("O1").EntireColumn.ColumnWidth = 16.7109375
$ws2.Range("P1").EntireColumn.ColumnWidth = 16.7109375

$ws2.Range("D8").Select()
$ws2.PageSetup.Orientation = 1

# --- Update sheet1 (SignIn) ---
$ws1 = $wb.Worksheets.Item("SignIn")

# Remove row 3
$ws1.Range("A3:D3").EntireRow.Delete()

# Update header D1 Passwords -> Password
$ws1.Range("D1").Value = "Password"

# Update row2 values
$ws1.Range("C2").Value = "testjaga006@gmail.com"
$ws1.Range("D2").Value = "jaga@12345"

# Re-add hyperlinks on row2 since values changed
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("C2"), "mailto:testjaga006@gmail.com")
$ws1.Hyperlinks.Add($ws1.Range("D2"), "mailto:jaga@12345")
$ws1.Range("C2").Style = "Hyperlink"
$ws1.Range("D2").Style = "Hyperlink"

$ws1.Range("A1").EntireColumn.ColumnWidth = 9.5703125
$ws1.Range("D1").EntireColumn.ColumnWidth = 11.42578125

$ws1.Range("D8").Select()

# --- Activate CreateAccount tab ---
$ws2.Activate()
